# Update cryptos list (price / 1h volume change) for the latest GitHub Actions run.
# A leading "'" forces a numeric-looking string (e.g. "168.00") to stay text,
# matching how the source sheet stores "Price" as literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.622.06'
$ws.Range("E2").Value = '  -2.68%  '
$ws.Range("D3").Value = '3.808.30'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D5").Value = '''600.20'
$ws.Range("E5").Value = '  -2.71%  '
$ws.Range("D6").Value = '''168.00'
$ws.Range("E6").Value = '  -5.35%  '
$ws.Range("D7").Value = '3.806.43'
$ws.Range("E7").Value = '  +0.70%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '''0.529'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  -4.89%  '
$ws.Range("E11").Value = '  -5.32%  '
$ws.Range("D12").Value = '''0.465'
$ws.Range("E12").Value = '  -4.25%  '
$ws.Range("D13").Value = '''38.37'
$ws.Range("E13").Value = '  -3.91%  '
$ws.Range("D14").Value = '''0.0000244'
$ws.Range("E14").Value = '  -4.42%  '
$ws.Range("D15").Value = '4.439.47'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '3.807.38'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '67.734.49'
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("D18").Value = '''7.23'
$ws.Range("E18").Value = '  -4.59%  '
$ws.Range("E19").Value = '  -3.78%  '
$ws.Range("D20").Value = '''17.37'
$ws.Range("E20").Value = '  +5.33%  '
$ws.Range("D21").Value = '''493.23'
$ws.Range("E21").Value = '  -3.25%  '
$ws.Range("D22").Value = '''9.44'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = '''0.738'
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("D24").Value = '''85.51'
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000147'
$ws.Range("E25").Value = '  +3.95%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '''2.38'
$ws.Range("E26").Value = '  -4.63%  '
$ws.Range("D27").Value = '''12.30'
$ws.Range("E27").Value = '  -4.73%  '
$ws.Range("D28").Value = '''10.11'
$ws.Range("E28").Value = '  -4.70%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").Value = '''2.98'
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("D31").Value = '''2.43'
$ws.Range("E31").Value = '  -4.12%  '
$ws.Range("D32").Value = '''32.78'
$ws.Range("E32").Value = '  +5.76%  '
$ws.Range("D33").Value = '''7.83'
$ws.Range("E33").Value = '  -3.08%  '
$ws.Range("D34").Value = '''0.109'
$ws.Range("E34").Value = '  -4.82%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  -3.01%  '
$ws.Range("D37").Value = '''5.82'
$ws.Range("E37").Value = '  -5.27%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.134'
$ws.Range("E38").Value = '  -5.05%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '''464.54'
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("D40").Value = '''0.329'
$ws.Range("E40").Value = '  -3.43%  '
$ws.Range("D41").Value = '''49.32'
$ws.Range("E41").Value = '  -1.00%  '
$ws.Range("D42").Value = '''2.00'
$ws.Range("E42").Value = '  -3.44%  '
$ws.Range("D43").Value = '''2.84'
$ws.Range("E43").Value = '  -4.85%  '
$ws.Range("D44").Value = '''8.39'
$ws.Range("E44").Value = '  -2.32%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = '''40.53'
$ws.Range("E46").Value = '  -8.71%  '
$ws.Range("D47").Value = '2.845.09'
$ws.Range("E47").Value = '  -3.84%  '
$ws.Range("D48").Value = '''139.97'
$ws.Range("E48").Value = '  +0.96%  '
$ws.Range("D49").Value = '''0.0350'
$ws.Range("E49").Value = '  -3.43%  '
$ws.Range("D50").Value = '''24.49'
$ws.Range("E50").Value = '  +12.16%  '
$ws.Range("D51").Value = '''25.72'
$ws.Range("E51").Value = '  -6.32%  '
